# Added 1.8V LDO to bom
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new BOM row (row 30) for the AAT3221 1.8V LDO regulator.
$ws.Range("B30").Value = "AAT3221"
$ws.Range("C30").Value = "SOT-23-5"
$ws.Range("F30").Value = "IC REG LDO 1.8V 0.15A SOT23-5"
$ws.Range("G30").Value = "Skyworks Solutions Inc"
$ws.Range("H30").Value = "AAT3221IGV-1.8-T1"
$ws.Range("I30").Value = "863-1508-1-ND"
$ws.Range("K30").Value = 0.17
$ws.Range("L30").Value = "-"
$ws.Range("M30").Value = 0.15
$ws.Range("N30").Value = 0.14

# Move/extend the active selection past the newly added row, like Excel
# does after typing data into the last row and pressing Enter/Down.
$null = $ws.Range("A31").Select()
